$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D values are written as text, not auto-converted numbers,
# and restore the original (unstyled) cell style afterward so formatting is unaffected.
$dRange = $ws.Range("D2:D51")
$dRange.NumberFormat = "@"

$ws.Range("D2").Value = '29.873.70'
$ws.Range("E2").Value = '  +0.01%  '

$ws.Range("D3").Value = '1.887.59'
$ws.Range("E3").Value = '  -0.32%  '

$ws.Range("E4").Value = '  +0.05%  '

$ws.Range("D5").Value = '0.7684'
$ws.Range("E5").Value = '  -0.87%  '

$ws.Range("D6").Value = '242.59'
$ws.Range("E6").Value = '  -0.92%  '

$ws.Range("D7").Value = '1.001'

$ws.Range("D8").Value = '0.3122'
$ws.Range("E8").Value = '  -0.87%  '

$ws.Range("D9").Value = '25.65'
$ws.Range("E9").Value = '  +0.54%  '

$ws.Range("D10").Value = '0.07163'
$ws.Range("E10").Value = '  -4.89%  '

$ws.Range("D11").Value = '0.08571'
$ws.Range("E11").Value = '  +5.63%  '

$ws.Range("D12").Value = '0.7631'
$ws.Range("E12").Value = '  -0.95%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '5.364'
$ws.Range("E13").Value = '  -2.05%  '

$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '1.877.04'
$ws.Range("E14").Value = '  -0.52%  '

$ws.Range("D15").Value = '93.58'
$ws.Range("E15").Value = '  +1.33%  '

$ws.Range("D16").Value = '6.150'
$ws.Range("E16").Value = '  -1.17%  '

$ws.Range("D17").Value = '29.834.87'
$ws.Range("E17").Value = '  +0.08%  '

$ws.Range("D18").Value = '13.77'
$ws.Range("E18").Value = '  -1.74%  '

$ws.Range("D19").Value = '244.42'
$ws.Range("E19").Value = '  -0.13%  '

$ws.Range("D20").Value = '0.000007799'
$ws.Range("E20").Value = '  -1.35%  '

$ws.Range("D21").Value = '2.141.24'
$ws.Range("E21").Value = '  +1.45%  '

$ws.Range("D22").Value = '0.9992'
$ws.Range("E22").Value = '  -0.05%  '

$ws.Range("D23").Value = '8.001'
$ws.Range("E23").Value = '  -1.24%  '

$ws.Range("E24").Value = '  +0.06%  '

$ws.Range("D25").Value = '0.1634'
$ws.Range("E25").Value = '  +4.09%  '

$ws.Range("D26").Value = '9.382'
$ws.Range("E26").Value = '  -0.81%  '

$ws.Range("D27").Value = '162.80'
$ws.Range("E27").Value = '  -0.07%  '

$ws.Range("E28").Value = '  -0.50%  '

$ws.Range("D29").Value = '2.030'
$ws.Range("E29").Value = '  -0.94%  '

$ws.Range("D30").Value = '1.462'
$ws.Range("E30").Value = '  +1.88%  '

$ws.Range("D31").Value = '1.539'
$ws.Range("E31").Value = '  -0.79%  '

$ws.Range("D32").Value = '4.511'
$ws.Range("E32").Value = '  +0.31%  '

$ws.Range("D33").Value = '4.092'
$ws.Range("E33").Value = '  -0.19%  '

$ws.Range("D34").Value = '0.05454'
$ws.Range("E34").Value = '  -1.09%  '

$ws.Range("E35").Value = '  -1.74%  '

$ws.Range("D36").Value = '0.7419'
$ws.Range("E36").Value = '  -2.09%  '

$ws.Range("D37").Value = '1.001'
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("D38").Value = '2.700'
$ws.Range("E38").Value = '  +2.25%  '

$ws.Range("D39").Value = '0.01954'
$ws.Range("E39").Value = '  +1.36%  '

$ws.Range("D40").Value = '2.782'
$ws.Range("E40").Value = '  -0.30%  '

$ws.Range("D41").Value = '0.4466'
$ws.Range("E41").Value = '  +0.23%  '

$ws.Range("D42").Value = '1.108.59'
$ws.Range("E42").Value = '  -4.66%  '

$ws.Range("B43").Value = 'Aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D43").Value = '72.98'
$ws.Range("E43").Value = '  -1.54%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = '6.069'
$ws.Range("E44").Value = '  +2.08%  '

$ws.Range("E46").Value = '  +0.03%  '

$ws.Range("D47").Value = '102.36'
$ws.Range("E47").Value = '  +0.06%  '

$ws.Range("D48").Value = '7.649'
$ws.Range("E48").Value = '  +1.45%  '

$ws.Range("D49").Value = '1.859'
$ws.Range("E49").Value = '  -2.45%  '

$ws.Range("D50").Value = '3.006'
$ws.Range("E50").Value = '  -3.60%  '

$ws.Range("D51").Value = '2.055.35'
$ws.Range("E51").Value = '  +1.49%  '

$dRange.Style = "Normal"
